# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets to
# reflect the refreshed data generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F7"  = 81
    "F8"  = 451
    "F13" = 295
    "F15" = 367
    "F18" = 6
    "F22" = 912
    "F25" = 323
    "F30" = 84
    "F32" = 243
    "F33" = 270
    "F34" = 1614
    "F39" = 299
    "F40" = 3603
    "F41" = 422
    "F43" = 905
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
